# Update master to output generated at 9a8706d
# Replace the three-digit x one-digit multiplication expressions in the
# practice worksheet table with newly generated problems.

$d = $word.ActiveDocument

$replacements = @(
    @("276×6=", "859×2="),
    @("447×3=", "661×7="),
    @("721×4=", "306×7="),
    @("533×7=", "489×4="),
    @("837×6=", "781×3="),
    @("155×2=", "502×2="),
    @("354×8=", "587×4="),
    @("942×6=", "359×3="),
    @("343×9=", "289×4="),
    @("960×2=", "195×5="),
    @("923×7=", "441×3="),
    @("915×9=", "399×9="),
    @("540×5=", "795×9="),
    @("273×5=", "991×7="),
    @("209×6=", "250×7="),
    @("301×4=", "966×4="),
    @("115×4=", "519×4="),
    @("701×8=", "163×5="),
    @("607×5=", "538×8="),
    @("494×5=", "819×7="),
    @("769×4=", "230×6="),
    @("898×6=", "133×2="),
    @("921×2=", "376×5="),
    @("661×9=", "493×8="),
    @("363×5=", "665×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
